$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.96000000000062
$ws.Range("G2").Value = 0.0006639630043193678
$ws.Range("H2").Value = 0.002146307009687996
$ws.Range("K2").Value = 4.453822739086347
$ws.Range("L2").Value = "[1.7469926267198037, 7.160652851452889]"
$ws.Range("M2").Value = 0.00133078966604705
$ws.Range("N2").Value = 0.00133078966604705
$ws.Range("O2").Value = -1.434000250287233
$ws.Range("P2").Value = "[-2.213895123250465, -0.6541053773240009]"
$ws.Range("Q2").Value = 0.0003442742813308808
$ws.Range("R2").Value = 0.0003442742813308808
$ws.Range("S2").Value = 13.71027461994294
$ws.Range("T2").Value = "[12.106905714240236, 15.313643525645647]"
$ws.Range("W2").Value = 5.924804804804946
$ws.Range("X2").Value = 2.702542542542607
$ws.Range("Y2").Value = 9.147067067067287

# Row 3 updates
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 22.80000000000013
$ws.Range("G3").Value = [double]"9.588033888285885e-06"
$ws.Range("H3").Value = 0.0001348797381065738
$ws.Range("I3").Value = 0.3740184267162057
$ws.Range("K3").Value = 6.260352112712463
$ws.Range("L3").Value = "[3.039134422295586, 9.48156980312934]"
$ws.Range("M3").Value = 0.0001604570184412513
$ws.Range("N3").Value = 0.0003209140368825025
$ws.Range("O3").Value = -1.245316006828387
$ws.Range("P3").Value = "[-1.8365266363327724, -0.6541053773240009]"
$ws.Range("Q3").Value = [double]"4.466460549501505e-05"
$ws.Range("R3").Value = [double]"8.932921099003011e-05"
$ws.Range("S3").Value = 12.32367413877413
$ws.Range("T3").Value = "[10.559271637210712, 14.08807664033755]"
$ws.Range("W3").Value = 4.518918918918946
$ws.Range("X3").Value = 2.373573573573588
$ws.Range("Y3").Value = 6.664264264264304
